$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 227 (shifts all existing rows 227..338 down to 228..339)
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new weekly data point.
$ws.Cells.Item(227, 1).Value = 4
$ws.Cells.Item(227, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(227, 3).Value = "Los Lagos"
$ws.Cells.Item(227, 4).Value = 44813
$ws.Cells.Item(227, 5).Value = 10
$ws.Cells.Item(227, 6).Value = 100112040
$ws.Cells.Item(227, 7).Value = "Cilantro"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 300
$ws.Cells.Item(227, 11).Value = 11000
$ws.Cells.Item(227, 12).Value = 12000
$ws.Cells.Item(227, 13).Value = 11500
$ws.Cells.Item(227, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(227, 15).Value = "Región Metropolitana"
$ws.Cells.Item(227, 16).Value = 319
$ws.Cells.Item(227, 17).Value = 36
$ws.Cells.Item(227, 18).Value = "Hortaliza"
